# Fruta / hortaliza, semanal
#
# A new weekly price report (fecha = 2022-09-07, serial 44811) is inserted
# as rows 16-18, ahead of the existing blocks, pushing every block that was
# at row 16 onward down by 3 rows (the old row 16 block is now at row 19,
# and so on, with one brand-new block appended at the very end of the
# sheet / the bottom of what used to be the last block).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at row 16. Excel shifts rows 16+ down to 19+; the new
# rows inherit formatting (incl. the column D date number format) from the
# row immediately above them, same as a manual Excel row insert.
$ws.Rows("16:18").Insert()

$newRows = @(
    @{ R = 16; Calidad = "Especial"; Volumen = 100; Pmin = 7000; Pmax = 8000; Pprom = 7500; Pkg = 2500 },
    @{ R = 17; Calidad = "Primera";  Volumen = 200; Pmin = 4000; Pmax = 5000; Pprom = 4500; Pkg = 1500 },
    @{ R = 18; Calidad = "Segunda";  Volumen = 200; Pmin = 3000; Pmax = 4000; Pprom = 3500; Pkg = 1167 }
)

foreach ($row in $newRows) {
    $r = $row.R
    $ws.Cells.Item($r, 1).Value = 1
    $ws.Cells.Item($r, 2).Value = "Agrícola del Norte S.A. de Arica"
    $ws.Cells.Item($r, 3).Value = "Arica y Parinacota"
    $ws.Cells.Item($r, 4).Value = 44811
    $ws.Cells.Item($r, 5).Value = 15
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100101
    $ws.Cells.Item($r, 8).Value = "Berries"
    $ws.Cells.Item($r, 9).Value = 100112025
    $ws.Cells.Item($r, 10).Value = "Frutilla"
    $ws.Cells.Item($r, 11).Value = "Sin especificar"
    $ws.Cells.Item($r, 12).Value = $row.Calidad
    $ws.Cells.Item($r, 13).Value = $row.Volumen
    $ws.Cells.Item($r, 14).Value = $row.Pmin
    $ws.Cells.Item($r, 15).Value = $row.Pmax
    $ws.Cells.Item($r, 16).Value = $row.Pprom
    $ws.Cells.Item($r, 17).Value = "`$/bandeja 3 kilos"
    $ws.Cells.Item($r, 18).Value = "Región de Arica y Parinacota"
    $ws.Cells.Item($r, 19).Value = $row.Pkg
    $ws.Cells.Item($r, 20).Value = 3
}

$ws.Range("D16:D18").NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Output "Inserted weekly block for 2022-09-07 at rows 16-18"
